# Applies cryptos list update (prices/volumes), matching commit:
# "Updated cryptos list on Tue Jan 23 03:44:25 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'40.058.10"
$ws.Range("E2").Value = "  -2.72%  "

$ws.Range("D3").Value = "'2.339.64"
$ws.Range("E3").Value = "  -3.88%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'310.85"
$ws.Range("E5").Value = "  -1.71%  "

$ws.Range("D6").Value = "'85.49"
$ws.Range("E6").Value = "  -3.58%  "

$ws.Range("E7").Value = "  -2.11%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("E9").Value = "  -2.27%  "

$ws.Range("D10").Value = "'0.0812"
$ws.Range("E10").Value = "  -2.43%  "

$ws.Range("D11").Value = "'30.02"
$ws.Range("E11").Value = "  -6.55%  "

$ws.Range("E12").Value = "  +1.07%  "

$ws.Range("D13").Value = "'2.696.96"
$ws.Range("E13").Value = "  -4.00%  "

$ws.Range("D14").Value = "'6.44"
$ws.Range("E14").Value = "  -4.17%  "

$ws.Range("D15").Value = "'14.79"
$ws.Range("E15").Value = "  -5.59%  "

$ws.Range("D16").Value = "'2.378.07"
$ws.Range("E16").Value = "  -2.50%  "

$ws.Range("D17").Value = "'0.758"
$ws.Range("E17").Value = "  -1.89%  "

$ws.Range("D18").Value = "'40.023.08"
$ws.Range("E18").Value = "  -2.66%  "

$ws.Range("E19").Value = "  -1.91%  "

$ws.Range("E20").Value = "  -1.75%  "

$ws.Range("D21").Value = "'67.95"
$ws.Range("E21").Value = "  -5.42%  "

$ws.Range("D22").Value = "'10.67"
$ws.Range("E22").Value = "  -3.18%  "

$ws.Range("D23").Value = "'235.28"
$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("E24").Value = "  -4.81%  "

$ws.Range("E25").Value = "  +0.10%  "

$ws.Range("E26").Value = "  -2.89%  "

$ws.Range("D27").Value = "'23.35"
$ws.Range("E27").Value = "  -2.56%  "

$ws.Range("E28").Value = "  -4.05%  "

$ws.Range("D29").Value = "'9.32"
$ws.Range("E29").Value = "  -2.46%  "

$ws.Range("D30").Value = "'34.90"
$ws.Range("E30").Value = "  +0.54%  "

$ws.Range("D31").Value = "'153.20"
$ws.Range("E31").Value = "  -2.38%  "

$ws.Range("E32").Value = "  -0.12%  "

$ws.Range("D33").Value = "'5.12"
$ws.Range("E33").Value = "  -2.64%  "

$ws.Range("E34").Value = "  -2.85%  "

$ws.Range("E35").Value = "  -3.23%  "

$ws.Range("E36").Value = "  -0.61%  "

$ws.Range("D37").Value = "'2.81"
$ws.Range("E37").Value = "  -3.71%  "

$ws.Range("D38").Value = "'0.0987"
$ws.Range("E38").Value = "  -1.25%  "

$ws.Range("D39").Value = "'15.64"
$ws.Range("E39").Value = "  -5.67%  "

$ws.Range("D40").Value = "'1.73"
$ws.Range("E40").Value = "  -2.59%  "

$ws.Range("D41").Value = "'3.90"
$ws.Range("E41").Value = "  +1.22%  "

$ws.Range("D42").Value = "'1.957.63"
$ws.Range("E42").Value = "  -1.28%  "

$ws.Range("E43").Value = "  -4.27%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'17.67"
$ws.Range("E44").Value = "  -3.34%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0264"
$ws.Range("E45").Value = "  -4.26%  "

$ws.Range("D46").Value = "'9.41"
$ws.Range("E46").Value = "  -0.79%  "

$ws.Range("D47").Value = "'2.71"
$ws.Range("E47").Value = "  -5.49%  "

$ws.Range("D48").Value = "'2.559.57"
$ws.Range("E48").Value = "  -4.15%  "

$ws.Range("D49").Value = "'92.86"
$ws.Range("E49").Value = "  -2.38%  "

$ws.Range("D50").Value = "'70.72"
$ws.Range("E50").Value = "  -3.18%  "

$ws.Range("D51").Value = "'51.05"
$ws.Range("E51").Value = "  -0.65%  "
